# Actualización automática 2025-08-19 11:50:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D10").Value = 1893.89
$wsGrupo.Range("D34").Value = "4 de 32"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F10").Value = 1893.89
$wsMensual.Range("F34").Value = 14129.6

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D3").Value = 4710.53
$wsCumpl.Range("E3").Value = -1590.4155
$wsCumpl.Range("F3").Value = 1.509729851260266

$wsCumpl.Range("D19").Value = 14249.68
$wsCumpl.Range("E19").Value = 17859.60107555788
$wsCumpl.Range("F19").Value = 0.4437869526405279
